$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I19").Value = 0.3532612570755742
$ws.Range("J19").Value = 0.2141232790717038
$ws.Range("K19").Value = -0.1106676956798791
$ws.Range("L19").Value = 2.89103650729813

$ws.Range("I20").Value = 0.6466138601298965
$ws.Range("J20").Value = 0.4806400749445675
$ws.Range("K20").Value = 0.2873446956070686
$ws.Range("L20").Value = 2.21951688451959
